$d = $word.ActiveDocument

# --- Change 1: "Projects and other experinces" -> "Other relevant experinces"
#     and move/insert the "_GoBack" bookmark right before this run.
$rng = $d.Content
$found = $rng.Find.Execute("Projects and other experinces", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $bmRange = $d.Range($rng.Start, $rng.Start)
    $d.Bookmarks.Add("_GoBack", $bmRange)
}
$d.Content.Find.Execute("Projects and other experinces", $true, $false, $false, $false, $false, $true, 1, $false, "Other relevant experinces", 2)

# --- Change 2: "Analyst | 09/2012 - 04/2013." -> "IT Analyst | 09/2012 - 04/2013."
$d.Content.Find.Execute("Analyst | 09/2012 - 04/2013.", $true, $false, $false, $false, $false, $true, 1, $false, "IT Analyst | 09/2012 - 04/2013.", 2)

# --- Change 3: merge "Network Administrator " + "| 01/2012 - 05/2012." runs
#     (which also removes the old "_GoBack" bookmark that sat between them)
#     into a single run reading "IT Analyst |Project 01/2012 - 05/2012."
$d.Content.Find.Execute("Network Administrator | 01/2012 - 05/2012.", $true, $false, $false, $false, $false, $true, 1, $false, "IT Analyst |Project 01/2012 - 05/2012.", 2)

# --- Change 4: mark the FollowedHyperlink style as a quick style (w:qFormat)
$followedHyperlinkStyle = $d.Styles.Item("FollowedHyperlink")
$followedHyperlinkStyle.QuickStyle = $true
